$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 335: update date, volumen, precio promedio ponderado and precio $/Kg
$ws.Range("D335").Value = 44505
$ws.Range("J335").Value = 1900
$ws.Range("M335").Value = 625
$ws.Range("P335").Value = 625

# Row 336: update date only
$ws.Range("D336").Value = 44505

# Rows 337-387: data shifted down by two rows (new entries inserted above),
# so rewrite each row fully with its new target content
# Row 337
$ws.Range("A337").Value = 3
$ws.Range("B337").Value = 'Femacal de La Calera'
$ws.Range("C337").Value = 'Coquimbo'
$ws.Range("D337").Value = 44487
$ws.Range("E337").Value = 5
$ws.Range("F337").Value = 100112006
$ws.Range("G337").Value = 'Repollo'
$ws.Range("H337").Value = 'Crespo record'
$ws.Range("I337").Value = 'Primera'
$ws.Range("J337").Value = 1830
$ws.Range("K337").Value = 600
$ws.Range("L337").Value = 650
$ws.Range("M337").Value = 624
$ws.Range("N337").Value = '$/unidad'
$ws.Range("O337").Value = 'Provincia de Quillota'
$ws.Range("P337").Value = 624
$ws.Range("Q337").Value = 1
$ws.Range("R337").Value = 'Hortaliza'

# Row 338
$ws.Range("A338").Value = 3
$ws.Range("B338").Value = 'Femacal de La Calera'
$ws.Range("C338").Value = 'Coquimbo'
$ws.Range("D338").Value = 44487
$ws.Range("E338").Value = 5
$ws.Range("F338").Value = 100112006
$ws.Range("G338").Value = 'Repollo'
$ws.Range("H338").Value = 'Crespo record'
$ws.Range("I338").Value = 'Segunda'
$ws.Range("J338").Value = 900
$ws.Range("K338").Value = 500
$ws.Range("L338").Value = 500
$ws.Range("M338").Value = 500
$ws.Range("N338").Value = '$/unidad'
$ws.Range("O338").Value = 'Provincia de Quillota'
$ws.Range("P338").Value = 500
$ws.Range("Q338").Value = 1
$ws.Range("R338").Value = 'Hortaliza'

# Row 339
$ws.Range("A339").Value = 3
$ws.Range("B339").Value = 'Femacal de La Calera'
$ws.Range("C339").Value = 'Coquimbo'
$ws.Range("D339").Value = 44425
$ws.Range("E339").Value = 5
$ws.Range("F339").Value = 100112006
$ws.Range("G339").Value = 'Repollo'
$ws.Range("H339").Value = 'Crespo record'
$ws.Range("I339").Value = 'Primera'
$ws.Range("J339").Value = 1800
$ws.Range("K339").Value = 700
$ws.Range("L339").Value = 700
$ws.Range("M339").Value = 700
$ws.Range("N339").Value = '$/unidad'
$ws.Range("O339").Value = 'Provincia de Quillota'
$ws.Range("P339").Value = 700
$ws.Range("Q339").Value = 1
$ws.Range("R339").Value = 'Hortaliza'

# Row 340
$ws.Range("A340").Value = 3
$ws.Range("B340").Value = 'Femacal de La Calera'
$ws.Range("C340").Value = 'Coquimbo'
$ws.Range("D340").Value = 44425
$ws.Range("E340").Value = 5
$ws.Range("F340").Value = 100112006
$ws.Range("G340").Value = 'Repollo'
$ws.Range("H340").Value = 'Crespo record'
$ws.Range("I340").Value = 'Segunda'
$ws.Range("J340").Value = 900
$ws.Range("K340").Value = 600
$ws.Range("L340").Value = 600
$ws.Range("M340").Value = 600
$ws.Range("N340").Value = '$/unidad'
$ws.Range("O340").Value = 'Provincia de Quillota'
$ws.Range("P340").Value = 600
$ws.Range("Q340").Value = 1
$ws.Range("R340").Value = 'Hortaliza'

# Row 341
$ws.Range("A341").Value = 3
$ws.Range("B341").Value = 'Femacal de La Calera'
$ws.Range("C341").Value = 'Coquimbo'
$ws.Range("D341").Value = 44343
$ws.Range("E341").Value = 5
$ws.Range("F341").Value = 100112006
$ws.Range("G341").Value = 'Repollo'
$ws.Range("H341").Value = 'Crespo record'
$ws.Range("I341").Value = 'Primera'
$ws.Range("J341").Value = 2500
$ws.Range("K341").Value = 800
$ws.Range("L341").Value = 800
$ws.Range("M341").Value = 800
$ws.Range("N341").Value = '$/unidad'
$ws.Range("O341").Value = 'Provincia de Quillota'
$ws.Range("P341").Value = 800
$ws.Range("Q341").Value = 1
$ws.Range("R341").Value = 'Hortaliza'

# Row 342
$ws.Range("A342").Value = 3
$ws.Range("B342").Value = 'Femacal de La Calera'
$ws.Range("C342").Value = 'Coquimbo'
$ws.Range("D342").Value = 44370
$ws.Range("E342").Value = 5
$ws.Range("F342").Value = 100112006
$ws.Range("G342").Value = 'Repollo'
$ws.Range("H342").Value = 'Crespo record'
$ws.Range("I342").Value = 'Primera'
$ws.Range("J342").Value = 1600
$ws.Range("K342").Value = 700
$ws.Range("L342").Value = 700
$ws.Range("M342").Value = 700
$ws.Range("N342").Value = '$/unidad'
$ws.Range("O342").Value = 'Provincia de Quillota'
$ws.Range("P342").Value = 700
$ws.Range("Q342").Value = 1
$ws.Range("R342").Value = 'Hortaliza'

# Row 343
$ws.Range("A343").Value = 3
$ws.Range("B343").Value = 'Femacal de La Calera'
$ws.Range("C343").Value = 'Coquimbo'
$ws.Range("D343").Value = 44370
$ws.Range("E343").Value = 5
$ws.Range("F343").Value = 100112006
$ws.Range("G343").Value = 'Repollo'
$ws.Range("H343").Value = 'Crespo record'
$ws.Range("I343").Value = 'Segunda'
$ws.Range("J343").Value = 1500
$ws.Range("K343").Value = 600
$ws.Range("L343").Value = 600
$ws.Range("M343").Value = 600
$ws.Range("N343").Value = '$/unidad'
$ws.Range("O343").Value = 'Provincia de Quillota'
$ws.Range("P343").Value = 600
$ws.Range("Q343").Value = 1
$ws.Range("R343").Value = 'Hortaliza'

# Row 344
$ws.Range("A344").Value = 3
$ws.Range("B344").Value = 'Femacal de La Calera'
$ws.Range("C344").Value = 'Coquimbo'
$ws.Range("D344").Value = 44449
$ws.Range("E344").Value = 5
$ws.Range("F344").Value = 100112006
$ws.Range("G344").Value = 'Repollo'
$ws.Range("H344").Value = 'Crespo record'
$ws.Range("I344").Value = 'Primera'
$ws.Range("J344").Value = 3400
$ws.Range("K344").Value = 500
$ws.Range("L344").Value = 600
$ws.Range("M344").Value = 553
$ws.Range("N344").Value = '$/unidad'
$ws.Range("O344").Value = 'Provincia de Quillota'
$ws.Range("P344").Value = 553
$ws.Range("Q344").Value = 1
$ws.Range("R344").Value = 'Hortaliza'

# Row 345
$ws.Range("A345").Value = 3
$ws.Range("B345").Value = 'Femacal de La Calera'
$ws.Range("C345").Value = 'Coquimbo'
$ws.Range("D345").Value = 44449
$ws.Range("E345").Value = 5
$ws.Range("F345").Value = 100112006
$ws.Range("G345").Value = 'Repollo'
$ws.Range("H345").Value = 'Crespo record'
$ws.Range("I345").Value = 'Segunda'
$ws.Range("J345").Value = 1200
$ws.Range("K345").Value = 400
$ws.Range("L345").Value = 400
$ws.Range("M345").Value = 400
$ws.Range("N345").Value = '$/unidad'
$ws.Range("O345").Value = 'Provincia de Quillota'
$ws.Range("P345").Value = 400
$ws.Range("Q345").Value = 1
$ws.Range("R345").Value = 'Hortaliza'

# Row 346
$ws.Range("A346").Value = 3
$ws.Range("B346").Value = 'Femacal de La Calera'
$ws.Range("C346").Value = 'Coquimbo'
$ws.Range("D346").Value = 44168
$ws.Range("E346").Value = 5
$ws.Range("F346").Value = 100112006
$ws.Range("G346").Value = 'Repollo'
$ws.Range("H346").Value = 'Crespo record'
$ws.Range("I346").Value = 'Primera'
$ws.Range("J346").Value = 1750
$ws.Range("K346").Value = 600
$ws.Range("L346").Value = 650
$ws.Range("M346").Value = 626
$ws.Range("N346").Value = '$/unidad'
$ws.Range("O346").Value = 'Provincia de Quillota'
$ws.Range("P346").Value = 626
$ws.Range("Q346").Value = 1
$ws.Range("R346").Value = 'Hortaliza'

# Row 347
$ws.Range("A347").Value = 3
$ws.Range("B347").Value = 'Femacal de La Calera'
$ws.Range("C347").Value = 'Coquimbo'
$ws.Range("D347").Value = 44168
$ws.Range("E347").Value = 5
$ws.Range("F347").Value = 100112006
$ws.Range("G347").Value = 'Repollo'
$ws.Range("H347").Value = 'Crespo record'
$ws.Range("I347").Value = 'Segunda'
$ws.Range("J347").Value = 900
$ws.Range("K347").Value = 550
$ws.Range("L347").Value = 550
$ws.Range("M347").Value = 550
$ws.Range("N347").Value = '$/unidad'
$ws.Range("O347").Value = 'Provincia de Quillota'
$ws.Range("P347").Value = 550
$ws.Range("Q347").Value = 1
$ws.Range("R347").Value = 'Hortaliza'

# Row 348
$ws.Range("A348").Value = 3
$ws.Range("B348").Value = 'Femacal de La Calera'
$ws.Range("C348").Value = 'Coquimbo'
$ws.Range("D348").Value = 44175
$ws.Range("E348").Value = 5
$ws.Range("F348").Value = 100112006
$ws.Range("G348").Value = 'Repollo'
$ws.Range("H348").Value = 'Crespo record'
$ws.Range("I348").Value = 'Primera'
$ws.Range("J348").Value = 1750
$ws.Range("K348").Value = 650
$ws.Range("L348").Value = 700
$ws.Range("M348").Value = 676
$ws.Range("N348").Value = '$/unidad'
$ws.Range("O348").Value = 'Provincia de Quillota'
$ws.Range("P348").Value = 676
$ws.Range("Q348").Value = 1
$ws.Range("R348").Value = 'Hortaliza'

# Row 349
$ws.Range("A349").Value = 3
$ws.Range("B349").Value = 'Femacal de La Calera'
$ws.Range("C349").Value = 'Coquimbo'
$ws.Range("D349").Value = 44392
$ws.Range("E349").Value = 5
$ws.Range("F349").Value = 100112006
$ws.Range("G349").Value = 'Repollo'
$ws.Range("H349").Value = 'Crespo record'
$ws.Range("I349").Value = 'Primera'
$ws.Range("J349").Value = 1300
$ws.Range("K349").Value = 600
$ws.Range("L349").Value = 600
$ws.Range("M349").Value = 600
$ws.Range("N349").Value = '$/unidad'
$ws.Range("O349").Value = 'Provincia de Quillota'
$ws.Range("P349").Value = 600
$ws.Range("Q349").Value = 1
$ws.Range("R349").Value = 'Hortaliza'

# Row 350
$ws.Range("A350").Value = 3
$ws.Range("B350").Value = 'Femacal de La Calera'
$ws.Range("C350").Value = 'Coquimbo'
$ws.Range("D350").Value = 44286
$ws.Range("E350").Value = 5
$ws.Range("F350").Value = 100112006
$ws.Range("G350").Value = 'Repollo'
$ws.Range("H350").Value = 'Crespo record'
$ws.Range("I350").Value = 'Primera'
$ws.Range("J350").Value = 900
$ws.Range("K350").Value = 900
$ws.Range("L350").Value = 900
$ws.Range("M350").Value = 900
$ws.Range("N350").Value = '$/unidad'
$ws.Range("O350").Value = 'Provincia de Quillota'
$ws.Range("P350").Value = 900
$ws.Range("Q350").Value = 1
$ws.Range("R350").Value = 'Hortaliza'

# Row 351
$ws.Range("A351").Value = 3
$ws.Range("B351").Value = 'Femacal de La Calera'
$ws.Range("C351").Value = 'Coquimbo'
$ws.Range("D351").Value = 44286
$ws.Range("E351").Value = 5
$ws.Range("F351").Value = 100112006
$ws.Range("G351").Value = 'Repollo'
$ws.Range("H351").Value = 'Crespo record'
$ws.Range("I351").Value = 'Segunda'
$ws.Range("J351").Value = 850
$ws.Range("K351").Value = 700
$ws.Range("L351").Value = 700
$ws.Range("M351").Value = 700
$ws.Range("N351").Value = '$/unidad'
$ws.Range("O351").Value = 'Provincia de Quillota'
$ws.Range("P351").Value = 700
$ws.Range("Q351").Value = 1
$ws.Range("R351").Value = 'Hortaliza'

# Row 352
$ws.Range("A352").Value = 3
$ws.Range("B352").Value = 'Femacal de La Calera'
$ws.Range("C352").Value = 'Coquimbo'
$ws.Range("D352").Value = 44473
$ws.Range("E352").Value = 5
$ws.Range("F352").Value = 100112006
$ws.Range("G352").Value = 'Repollo'
$ws.Range("H352").Value = 'Crespo record'
$ws.Range("I352").Value = 'Primera'
$ws.Range("J352").Value = 1300
$ws.Range("K352").Value = 600
$ws.Range("L352").Value = 600
$ws.Range("M352").Value = 600
$ws.Range("N352").Value = '$/unidad'
$ws.Range("O352").Value = 'Provincia de Quillota'
$ws.Range("P352").Value = 600
$ws.Range("Q352").Value = 1
$ws.Range("R352").Value = 'Hortaliza'

# Row 353
$ws.Range("A353").Value = 3
$ws.Range("B353").Value = 'Femacal de La Calera'
$ws.Range("C353").Value = 'Coquimbo'
$ws.Range("D353").Value = 44473
$ws.Range("E353").Value = 5
$ws.Range("F353").Value = 100112006
$ws.Range("G353").Value = 'Repollo'
$ws.Range("H353").Value = 'Crespo record'
$ws.Range("I353").Value = 'Segunda'
$ws.Range("J353").Value = 1200
$ws.Range("K353").Value = 500
$ws.Range("L353").Value = 500
$ws.Range("M353").Value = 500
$ws.Range("N353").Value = '$/unidad'
$ws.Range("O353").Value = 'Provincia de Quillota'
$ws.Range("P353").Value = 500
$ws.Range("Q353").Value = 1
$ws.Range("R353").Value = 'Hortaliza'

# Row 354
$ws.Range("A354").Value = 3
$ws.Range("B354").Value = 'Femacal de La Calera'
$ws.Range("C354").Value = 'Coquimbo'
$ws.Range("D354").Value = 44400
$ws.Range("E354").Value = 5
$ws.Range("F354").Value = 100112006
$ws.Range("G354").Value = 'Repollo'
$ws.Range("H354").Value = 'Crespo record'
$ws.Range("I354").Value = 'Primera'
$ws.Range("J354").Value = 1850
$ws.Range("K354").Value = 700
$ws.Range("L354").Value = 750
$ws.Range("M354").Value = 724
$ws.Range("N354").Value = '$/unidad'
$ws.Range("O354").Value = 'Provincia de Quillota'
$ws.Range("P354").Value = 724
$ws.Range("Q354").Value = 1
$ws.Range("R354").Value = 'Hortaliza'

# Row 355
$ws.Range("A355").Value = 3
$ws.Range("B355").Value = 'Femacal de La Calera'
$ws.Range("C355").Value = 'Coquimbo'
$ws.Range("D355").Value = 44484
$ws.Range("E355").Value = 5
$ws.Range("F355").Value = 100112006
$ws.Range("G355").Value = 'Repollo'
$ws.Range("H355").Value = 'Crespo record'
$ws.Range("I355").Value = 'Primera'
$ws.Range("J355").Value = 1800
$ws.Range("K355").Value = 600
$ws.Range("L355").Value = 650
$ws.Range("M355").Value = 624
$ws.Range("N355").Value = '$/unidad'
$ws.Range("O355").Value = 'Provincia de Quillota'
$ws.Range("P355").Value = 624
$ws.Range("Q355").Value = 1
$ws.Range("R355").Value = 'Hortaliza'

# Row 356
$ws.Range("A356").Value = 3
$ws.Range("B356").Value = 'Femacal de La Calera'
$ws.Range("C356").Value = 'Coquimbo'
$ws.Range("D356").Value = 44484
$ws.Range("E356").Value = 5
$ws.Range("F356").Value = 100112006
$ws.Range("G356").Value = 'Repollo'
$ws.Range("H356").Value = 'Crespo record'
$ws.Range("I356").Value = 'Segunda'
$ws.Range("J356").Value = 900
$ws.Range("K356").Value = 500
$ws.Range("L356").Value = 500
$ws.Range("M356").Value = 500
$ws.Range("N356").Value = '$/unidad'
$ws.Range("O356").Value = 'Provincia de Quillota'
$ws.Range("P356").Value = 500
$ws.Range("Q356").Value = 1
$ws.Range("R356").Value = 'Hortaliza'

# Row 357
$ws.Range("A357").Value = 3
$ws.Range("B357").Value = 'Femacal de La Calera'
$ws.Range("C357").Value = 'Coquimbo'
$ws.Range("D357").Value = 44181
$ws.Range("E357").Value = 5
$ws.Range("F357").Value = 100112006
$ws.Range("G357").Value = 'Repollo'
$ws.Range("H357").Value = 'Crespo record'
$ws.Range("I357").Value = 'Primera'
$ws.Range("J357").Value = 1750
$ws.Range("K357").Value = 600
$ws.Range("L357").Value = 650
$ws.Range("M357").Value = 626
$ws.Range("N357").Value = '$/unidad'
$ws.Range("O357").Value = 'Provincia de Quillota'
$ws.Range("P357").Value = 626
$ws.Range("Q357").Value = 1
$ws.Range("R357").Value = 'Hortaliza'

# Row 358
$ws.Range("A358").Value = 3
$ws.Range("B358").Value = 'Femacal de La Calera'
$ws.Range("C358").Value = 'Coquimbo'
$ws.Range("D358").Value = 44181
$ws.Range("E358").Value = 5
$ws.Range("F358").Value = 100112006
$ws.Range("G358").Value = 'Repollo'
$ws.Range("H358").Value = 'Crespo record'
$ws.Range("I358").Value = 'Segunda'
$ws.Range("J358").Value = 900
$ws.Range("K358").Value = 550
$ws.Range("L358").Value = 550
$ws.Range("M358").Value = 550
$ws.Range("N358").Value = '$/unidad'
$ws.Range("O358").Value = 'Provincia de Quillota'
$ws.Range("P358").Value = 550
$ws.Range("Q358").Value = 1
$ws.Range("R358").Value = 'Hortaliza'

# Row 359
$ws.Range("A359").Value = 3
$ws.Range("B359").Value = 'Femacal de La Calera'
$ws.Range("C359").Value = 'Coquimbo'
$ws.Range("D359").Value = 44494
$ws.Range("E359").Value = 5
$ws.Range("F359").Value = 100112006
$ws.Range("G359").Value = 'Repollo'
$ws.Range("H359").Value = 'Crespo record'
$ws.Range("I359").Value = 'Primera'
$ws.Range("J359").Value = 1800
$ws.Range("K359").Value = 600
$ws.Range("L359").Value = 600
$ws.Range("M359").Value = 600
$ws.Range("N359").Value = '$/unidad'
$ws.Range("O359").Value = 'Provincia de Quillota'
$ws.Range("P359").Value = 600
$ws.Range("Q359").Value = 1
$ws.Range("R359").Value = 'Hortaliza'

# Row 360
$ws.Range("A360").Value = 3
$ws.Range("B360").Value = 'Femacal de La Calera'
$ws.Range("C360").Value = 'Coquimbo'
$ws.Range("D360").Value = 44494
$ws.Range("E360").Value = 5
$ws.Range("F360").Value = 100112006
$ws.Range("G360").Value = 'Repollo'
$ws.Range("H360").Value = 'Crespo record'
$ws.Range("I360").Value = 'Segunda'
$ws.Range("J360").Value = 1200
$ws.Range("K360").Value = 500
$ws.Range("L360").Value = 500
$ws.Range("M360").Value = 500
$ws.Range("N360").Value = '$/unidad'
$ws.Range("O360").Value = 'Provincia de Quillota'
$ws.Range("P360").Value = 500
$ws.Range("Q360").Value = 1
$ws.Range("R360").Value = 'Hortaliza'

# Row 361
$ws.Range("A361").Value = 3
$ws.Range("B361").Value = 'Femacal de La Calera'
$ws.Range("C361").Value = 'Coquimbo'
$ws.Range("D361").Value = 44342
$ws.Range("E361").Value = 5
$ws.Range("F361").Value = 100112006
$ws.Range("G361").Value = 'Repollo'
$ws.Range("H361").Value = 'Crespo record'
$ws.Range("I361").Value = 'Primera'
$ws.Range("J361").Value = 1300
$ws.Range("K361").Value = 700
$ws.Range("L361").Value = 700
$ws.Range("M361").Value = 700
$ws.Range("N361").Value = '$/unidad'
$ws.Range("O361").Value = 'Provincia de Quillota'
$ws.Range("P361").Value = 700
$ws.Range("Q361").Value = 1
$ws.Range("R361").Value = 'Hortaliza'

# Row 362
$ws.Range("A362").Value = 3
$ws.Range("B362").Value = 'Femacal de La Calera'
$ws.Range("C362").Value = 'Coquimbo'
$ws.Range("D362").Value = 44445
$ws.Range("E362").Value = 5
$ws.Range("F362").Value = 100112006
$ws.Range("G362").Value = 'Repollo'
$ws.Range("H362").Value = 'Crespo record'
$ws.Range("I362").Value = 'Primera'
$ws.Range("J362").Value = 1850
$ws.Range("K362").Value = 600
$ws.Range("L362").Value = 600
$ws.Range("M362").Value = 600
$ws.Range("N362").Value = '$/unidad'
$ws.Range("O362").Value = 'Provincia de Quillota'
$ws.Range("P362").Value = 600
$ws.Range("Q362").Value = 1
$ws.Range("R362").Value = 'Hortaliza'

# Row 363
$ws.Range("A363").Value = 3
$ws.Range("B363").Value = 'Femacal de La Calera'
$ws.Range("C363").Value = 'Coquimbo'
$ws.Range("D363").Value = 44445
$ws.Range("E363").Value = 5
$ws.Range("F363").Value = 100112006
$ws.Range("G363").Value = 'Repollo'
$ws.Range("H363").Value = 'Crespo record'
$ws.Range("I363").Value = 'Segunda'
$ws.Range("J363").Value = 1900
$ws.Range("K363").Value = 500
$ws.Range("L363").Value = 500
$ws.Range("M363").Value = 500
$ws.Range("N363").Value = '$/unidad'
$ws.Range("O363").Value = 'Provincia de Quillota'
$ws.Range("P363").Value = 500
$ws.Range("Q363").Value = 1
$ws.Range("R363").Value = 'Hortaliza'

# Row 364
$ws.Range("A364").Value = 3
$ws.Range("B364").Value = 'Femacal de La Calera'
$ws.Range("C364").Value = 'Coquimbo'
$ws.Range("D364").Value = 44328
$ws.Range("E364").Value = 5
$ws.Range("F364").Value = 100112006
$ws.Range("G364").Value = 'Repollo'
$ws.Range("H364").Value = 'Crespo record'
$ws.Range("I364").Value = 'Primera'
$ws.Range("J364").Value = 1800
$ws.Range("K364").Value = 800
$ws.Range("L364").Value = 800
$ws.Range("M364").Value = 800
$ws.Range("N364").Value = '$/unidad'
$ws.Range("O364").Value = 'Provincia de Quillota'
$ws.Range("P364").Value = 800
$ws.Range("Q364").Value = 1
$ws.Range("R364").Value = 'Hortaliza'

# Row 365
$ws.Range("A365").Value = 3
$ws.Range("B365").Value = 'Femacal de La Calera'
$ws.Range("C365").Value = 'Coquimbo'
$ws.Range("D365").Value = 44301
$ws.Range("E365").Value = 5
$ws.Range("F365").Value = 100112006
$ws.Range("G365").Value = 'Repollo'
$ws.Range("H365").Value = 'Crespo record'
$ws.Range("I365").Value = 'Primera'
$ws.Range("J365").Value = 1200
$ws.Range("K365").Value = 900
$ws.Range("L365").Value = 900
$ws.Range("M365").Value = 900
$ws.Range("N365").Value = '$/unidad'
$ws.Range("O365").Value = 'Provincia de Quillota'
$ws.Range("P365").Value = 900
$ws.Range("Q365").Value = 1
$ws.Range("R365").Value = 'Hortaliza'

# Row 366
$ws.Range("A366").Value = 3
$ws.Range("B366").Value = 'Femacal de La Calera'
$ws.Range("C366").Value = 'Coquimbo'
$ws.Range("D366").Value = 44301
$ws.Range("E366").Value = 5
$ws.Range("F366").Value = 100112006
$ws.Range("G366").Value = 'Repollo'
$ws.Range("H366").Value = 'Crespo record'
$ws.Range("I366").Value = 'Segunda'
$ws.Range("J366").Value = 950
$ws.Range("K366").Value = 700
$ws.Range("L366").Value = 700
$ws.Range("M366").Value = 700
$ws.Range("N366").Value = '$/unidad'
$ws.Range("O366").Value = 'Provincia de Quillota'
$ws.Range("P366").Value = 700
$ws.Range("Q366").Value = 1
$ws.Range("R366").Value = 'Hortaliza'

# Row 367
$ws.Range("A367").Value = 3
$ws.Range("B367").Value = 'Femacal de La Calera'
$ws.Range("C367").Value = 'Coquimbo'
$ws.Range("D367").Value = 44330
$ws.Range("E367").Value = 5
$ws.Range("F367").Value = 100112006
$ws.Range("G367").Value = 'Repollo'
$ws.Range("H367").Value = 'Crespo record'
$ws.Range("I367").Value = 'Primera'
$ws.Range("J367").Value = 1300
$ws.Range("K367").Value = 800
$ws.Range("L367").Value = 800
$ws.Range("M367").Value = 800
$ws.Range("N367").Value = '$/unidad'
$ws.Range("O367").Value = 'Provincia de Quillota'
$ws.Range("P367").Value = 800
$ws.Range("Q367").Value = 1
$ws.Range("R367").Value = 'Hortaliza'

# Row 368
$ws.Range("A368").Value = 3
$ws.Range("B368").Value = 'Femacal de La Calera'
$ws.Range("C368").Value = 'Coquimbo'
$ws.Range("D368").Value = 44330
$ws.Range("E368").Value = 5
$ws.Range("F368").Value = 100112006
$ws.Range("G368").Value = 'Repollo'
$ws.Range("H368").Value = 'Crespo record'
$ws.Range("I368").Value = 'Segunda'
$ws.Range("J368").Value = 1200
$ws.Range("K368").Value = 650
$ws.Range("L368").Value = 650
$ws.Range("M368").Value = 650
$ws.Range("N368").Value = '$/unidad'
$ws.Range("O368").Value = 'Provincia de Quillota'
$ws.Range("P368").Value = 650
$ws.Range("Q368").Value = 1
$ws.Range("R368").Value = 'Hortaliza'

# Row 369
$ws.Range("A369").Value = 3
$ws.Range("B369").Value = 'Femacal de La Calera'
$ws.Range("C369").Value = 'Coquimbo'
$ws.Range("D369").Value = 44270
$ws.Range("E369").Value = 5
$ws.Range("F369").Value = 100112006
$ws.Range("G369").Value = 'Repollo'
$ws.Range("H369").Value = 'Crespo record'
$ws.Range("I369").Value = 'Primera'
$ws.Range("J369").Value = 3400
$ws.Range("K369").Value = 800
$ws.Range("L369").Value = 850
$ws.Range("M369").Value = 824
$ws.Range("N369").Value = '$/unidad'
$ws.Range("O369").Value = 'Provincia de Quillota'
$ws.Range("P369").Value = 824
$ws.Range("Q369").Value = 1
$ws.Range("R369").Value = 'Hortaliza'

# Row 370
$ws.Range("A370").Value = 3
$ws.Range("B370").Value = 'Femacal de La Calera'
$ws.Range("C370").Value = 'Coquimbo'
$ws.Range("D370").Value = 44295
$ws.Range("E370").Value = 5
$ws.Range("F370").Value = 100112006
$ws.Range("G370").Value = 'Repollo'
$ws.Range("H370").Value = 'Crespo record'
$ws.Range("I370").Value = 'Primera'
$ws.Range("J370").Value = 1600
$ws.Range("K370").Value = 800
$ws.Range("L370").Value = 800
$ws.Range("M370").Value = 800
$ws.Range("N370").Value = '$/unidad'
$ws.Range("O370").Value = 'Provincia de Quillota'
$ws.Range("P370").Value = 800
$ws.Range("Q370").Value = 1
$ws.Range("R370").Value = 'Hortaliza'

# Row 371
$ws.Range("A371").Value = 3
$ws.Range("B371").Value = 'Femacal de La Calera'
$ws.Range("C371").Value = 'Coquimbo'
$ws.Range("D371").Value = 44217
$ws.Range("E371").Value = 5
$ws.Range("F371").Value = 100112006
$ws.Range("G371").Value = 'Repollo'
$ws.Range("H371").Value = 'Crespo record'
$ws.Range("I371").Value = 'Primera'
$ws.Range("J371").Value = 1200
$ws.Range("K371").Value = 800
$ws.Range("L371").Value = 800
$ws.Range("M371").Value = 800
$ws.Range("N371").Value = '$/unidad'
$ws.Range("O371").Value = 'Provincia de Quillota'
$ws.Range("P371").Value = 800
$ws.Range("Q371").Value = 1
$ws.Range("R371").Value = 'Hortaliza'

# Row 372
$ws.Range("A372").Value = 3
$ws.Range("B372").Value = 'Femacal de La Calera'
$ws.Range("C372").Value = 'Coquimbo'
$ws.Range("D372").Value = 44217
$ws.Range("E372").Value = 5
$ws.Range("F372").Value = 100112006
$ws.Range("G372").Value = 'Repollo'
$ws.Range("H372").Value = 'Crespo record'
$ws.Range("I372").Value = 'Segunda'
$ws.Range("J372").Value = 950
$ws.Range("K372").Value = 700
$ws.Range("L372").Value = 700
$ws.Range("M372").Value = 700
$ws.Range("N372").Value = '$/unidad'
$ws.Range("O372").Value = 'Provincia de Quillota'
$ws.Range("P372").Value = 700
$ws.Range("Q372").Value = 1
$ws.Range("R372").Value = 'Hortaliza'

# Row 373
$ws.Range("A373").Value = 3
$ws.Range("B373").Value = 'Femacal de La Calera'
$ws.Range("C373").Value = 'Coquimbo'
$ws.Range("D373").Value = 44421
$ws.Range("E373").Value = 5
$ws.Range("F373").Value = 100112006
$ws.Range("G373").Value = 'Repollo'
$ws.Range("H373").Value = 'Crespo record'
$ws.Range("I373").Value = 'Primera'
$ws.Range("J373").Value = 1600
$ws.Range("K373").Value = 700
$ws.Range("L373").Value = 700
$ws.Range("M373").Value = 700
$ws.Range("N373").Value = '$/unidad'
$ws.Range("O373").Value = 'Provincia de Quillota'
$ws.Range("P373").Value = 700
$ws.Range("Q373").Value = 1
$ws.Range("R373").Value = 'Hortaliza'

# Row 374
$ws.Range("A374").Value = 3
$ws.Range("B374").Value = 'Femacal de La Calera'
$ws.Range("C374").Value = 'Coquimbo'
$ws.Range("D374").Value = 44383
$ws.Range("E374").Value = 5
$ws.Range("F374").Value = 100112006
$ws.Range("G374").Value = 'Repollo'
$ws.Range("H374").Value = 'Crespo record'
$ws.Range("I374").Value = 'Primera'
$ws.Range("J374").Value = 1400
$ws.Range("K374").Value = 700
$ws.Range("L374").Value = 700
$ws.Range("M374").Value = 700
$ws.Range("N374").Value = '$/unidad'
$ws.Range("O374").Value = 'Provincia de Quillota'
$ws.Range("P374").Value = 700
$ws.Range("Q374").Value = 1
$ws.Range("R374").Value = 'Hortaliza'

# Row 375
$ws.Range("A375").Value = 3
$ws.Range("B375").Value = 'Femacal de La Calera'
$ws.Range("C375").Value = 'Coquimbo'
$ws.Range("D375").Value = 44383
$ws.Range("E375").Value = 5
$ws.Range("F375").Value = 100112006
$ws.Range("G375").Value = 'Repollo'
$ws.Range("H375").Value = 'Crespo record'
$ws.Range("I375").Value = 'Segunda'
$ws.Range("J375").Value = 1250
$ws.Range("K375").Value = 600
$ws.Range("L375").Value = 600
$ws.Range("M375").Value = 600
$ws.Range("N375").Value = '$/unidad'
$ws.Range("O375").Value = 'Provincia de Quillota'
$ws.Range("P375").Value = 600
$ws.Range("Q375").Value = 1
$ws.Range("R375").Value = 'Hortaliza'

# Row 376
$ws.Range("A376").Value = 3
$ws.Range("B376").Value = 'Femacal de La Calera'
$ws.Range("C376").Value = 'Coquimbo'
$ws.Range("D376").Value = 44244
$ws.Range("E376").Value = 5
$ws.Range("F376").Value = 100112006
$ws.Range("G376").Value = 'Repollo'
$ws.Range("H376").Value = 'Crespo record'
$ws.Range("I376").Value = 'Primera'
$ws.Range("J376").Value = 950
$ws.Range("K376").Value = 800
$ws.Range("L376").Value = 800
$ws.Range("M376").Value = 800
$ws.Range("N376").Value = '$/unidad'
$ws.Range("O376").Value = 'Provincia de Quillota'
$ws.Range("P376").Value = 800
$ws.Range("Q376").Value = 1
$ws.Range("R376").Value = 'Hortaliza'

# Row 377
$ws.Range("A377").Value = 3
$ws.Range("B377").Value = 'Femacal de La Calera'
$ws.Range("C377").Value = 'Coquimbo'
$ws.Range("D377").Value = 44244
$ws.Range("E377").Value = 5
$ws.Range("F377").Value = 100112006
$ws.Range("G377").Value = 'Repollo'
$ws.Range("H377").Value = 'Crespo record'
$ws.Range("I377").Value = 'Segunda'
$ws.Range("J377").Value = 850
$ws.Range("K377").Value = 700
$ws.Range("L377").Value = 700
$ws.Range("M377").Value = 700
$ws.Range("N377").Value = '$/unidad'
$ws.Range("O377").Value = 'Provincia de Quillota'
$ws.Range("P377").Value = 700
$ws.Range("Q377").Value = 1
$ws.Range("R377").Value = 'Hortaliza'

# Row 378
$ws.Range("A378").Value = 3
$ws.Range("B378").Value = 'Femacal de La Calera'
$ws.Range("C378").Value = 'Coquimbo'
$ws.Range("D378").Value = 44307
$ws.Range("E378").Value = 5
$ws.Range("F378").Value = 100112006
$ws.Range("G378").Value = 'Repollo'
$ws.Range("H378").Value = 'Crespo record'
$ws.Range("I378").Value = 'Primera'
$ws.Range("J378").Value = 1300
$ws.Range("K378").Value = 900
$ws.Range("L378").Value = 900
$ws.Range("M378").Value = 900
$ws.Range("N378").Value = '$/unidad'
$ws.Range("O378").Value = 'Provincia de Quillota'
$ws.Range("P378").Value = 900
$ws.Range("Q378").Value = 1
$ws.Range("R378").Value = 'Hortaliza'

# Row 379
$ws.Range("A379").Value = 3
$ws.Range("B379").Value = 'Femacal de La Calera'
$ws.Range("C379").Value = 'Coquimbo'
$ws.Range("D379").Value = 44307
$ws.Range("E379").Value = 5
$ws.Range("F379").Value = 100112006
$ws.Range("G379").Value = 'Repollo'
$ws.Range("H379").Value = 'Crespo record'
$ws.Range("I379").Value = 'Segunda'
$ws.Range("J379").Value = 850
$ws.Range("K379").Value = 700
$ws.Range("L379").Value = 700
$ws.Range("M379").Value = 700
$ws.Range("N379").Value = '$/unidad'
$ws.Range("O379").Value = 'Provincia de Quillota'
$ws.Range("P379").Value = 700
$ws.Range("Q379").Value = 1
$ws.Range("R379").Value = 'Hortaliza'

# Row 380
$ws.Range("A380").Value = 3
$ws.Range("B380").Value = 'Femacal de La Calera'
$ws.Range("C380").Value = 'Coquimbo'
$ws.Range("D380").Value = 44273
$ws.Range("E380").Value = 5
$ws.Range("F380").Value = 100112006
$ws.Range("G380").Value = 'Repollo'
$ws.Range("H380").Value = 'Crespo record'
$ws.Range("I380").Value = 'Primera'
$ws.Range("J380").Value = 2400
$ws.Range("K380").Value = 850
$ws.Range("L380").Value = 900
$ws.Range("M380").Value = 875
$ws.Range("N380").Value = '$/unidad'
$ws.Range("O380").Value = 'Provincia de Quillota'
$ws.Range("P380").Value = 875
$ws.Range("Q380").Value = 1
$ws.Range("R380").Value = 'Hortaliza'

# Row 381
$ws.Range("A381").Value = 3
$ws.Range("B381").Value = 'Femacal de La Calera'
$ws.Range("C381").Value = 'Coquimbo'
$ws.Range("D381").Value = 44273
$ws.Range("E381").Value = 5
$ws.Range("F381").Value = 100112006
$ws.Range("G381").Value = 'Repollo'
$ws.Range("H381").Value = 'Crespo record'
$ws.Range("I381").Value = 'Segunda'
$ws.Range("J381").Value = 900
$ws.Range("K381").Value = 700
$ws.Range("L381").Value = 700
$ws.Range("M381").Value = 700
$ws.Range("N381").Value = '$/unidad'
$ws.Range("O381").Value = 'Provincia de Quillota'
$ws.Range("P381").Value = 700
$ws.Range("Q381").Value = 1
$ws.Range("R381").Value = 'Hortaliza'

# Row 382
$ws.Range("A382").Value = 3
$ws.Range("B382").Value = 'Femacal de La Calera'
$ws.Range("C382").Value = 'Coquimbo'
$ws.Range("D382").Value = 44433
$ws.Range("E382").Value = 5
$ws.Range("F382").Value = 100112006
$ws.Range("G382").Value = 'Repollo'
$ws.Range("H382").Value = 'Crespo record'
$ws.Range("I382").Value = 'Primera'
$ws.Range("J382").Value = 1200
$ws.Range("K382").Value = 700
$ws.Range("L382").Value = 700
$ws.Range("M382").Value = 700
$ws.Range("N382").Value = '$/unidad'
$ws.Range("O382").Value = 'Provincia de Quillota'
$ws.Range("P382").Value = 700
$ws.Range("Q382").Value = 1
$ws.Range("R382").Value = 'Hortaliza'

# Row 383
$ws.Range("A383").Value = 3
$ws.Range("B383").Value = 'Femacal de La Calera'
$ws.Range("C383").Value = 'Coquimbo'
$ws.Range("D383").Value = 44433
$ws.Range("E383").Value = 5
$ws.Range("F383").Value = 100112006
$ws.Range("G383").Value = 'Repollo'
$ws.Range("H383").Value = 'Crespo record'
$ws.Range("I383").Value = 'Segunda'
$ws.Range("J383").Value = 900
$ws.Range("K383").Value = 500
$ws.Range("L383").Value = 500
$ws.Range("M383").Value = 500
$ws.Range("N383").Value = '$/unidad'
$ws.Range("O383").Value = 'Provincia de Quillota'
$ws.Range("P383").Value = 500
$ws.Range("Q383").Value = 1
$ws.Range("R383").Value = 'Hortaliza'

# Row 384
$ws.Range("A384").Value = 3
$ws.Range("B384").Value = 'Femacal de La Calera'
$ws.Range("C384").Value = 'Coquimbo'
$ws.Range("D384").Value = 44302
$ws.Range("E384").Value = 5
$ws.Range("F384").Value = 100112006
$ws.Range("G384").Value = 'Repollo'
$ws.Range("H384").Value = 'Crespo record'
$ws.Range("I384").Value = 'Primera'
$ws.Range("J384").Value = 1300
$ws.Range("K384").Value = 900
$ws.Range("L384").Value = 900
$ws.Range("M384").Value = 900
$ws.Range("N384").Value = '$/unidad'
$ws.Range("O384").Value = 'Provincia de Quillota'
$ws.Range("P384").Value = 900
$ws.Range("Q384").Value = 1
$ws.Range("R384").Value = 'Hortaliza'

# Row 385
$ws.Range("A385").Value = 3
$ws.Range("B385").Value = 'Femacal de La Calera'
$ws.Range("C385").Value = 'Coquimbo'
$ws.Range("D385").Value = 44179
$ws.Range("E385").Value = 5
$ws.Range("F385").Value = 100112006
$ws.Range("G385").Value = 'Repollo'
$ws.Range("H385").Value = 'Crespo record'
$ws.Range("I385").Value = 'Primera'
$ws.Range("J385").Value = 1600
$ws.Range("K385").Value = 600
$ws.Range("L385").Value = 700
$ws.Range("M385").Value = 653
$ws.Range("N385").Value = '$/unidad'
$ws.Range("O385").Value = 'Provincia de Quillota'
$ws.Range("P385").Value = 653
$ws.Range("Q385").Value = 1
$ws.Range("R385").Value = 'Hortaliza'

# Row 386
$ws.Range("A386").Value = 3
$ws.Range("B386").Value = 'Femacal de La Calera'
$ws.Range("C386").Value = 'Coquimbo'
$ws.Range("D386").Value = 44179
$ws.Range("D386").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E386").Value = 5
$ws.Range("F386").Value = 100112006
$ws.Range("G386").Value = 'Repollo'
$ws.Range("H386").Value = 'Crespo record'
$ws.Range("I386").Value = 'Segunda'
$ws.Range("J386").Value = 1690
$ws.Range("K386").Value = 500
$ws.Range("L386").Value = 550
$ws.Range("M386").Value = 529
$ws.Range("N386").Value = '$/unidad'
$ws.Range("O386").Value = 'Provincia de Quillota'
$ws.Range("P386").Value = 529
$ws.Range("Q386").Value = 1
$ws.Range("R386").Value = 'Hortaliza'

# Row 387
$ws.Range("A387").Value = 3
$ws.Range("B387").Value = 'Femacal de La Calera'
$ws.Range("C387").Value = 'Coquimbo'
$ws.Range("D387").Value = 44491
$ws.Range("D387").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E387").Value = 5
$ws.Range("F387").Value = 100112006
$ws.Range("G387").Value = 'Repollo'
$ws.Range("H387").Value = 'Crespo record'
$ws.Range("I387").Value = 'Primera'
$ws.Range("J387").Value = 4650
$ws.Range("K387").Value = 500
$ws.Range("L387").Value = 650
$ws.Range("M387").Value = 585
$ws.Range("N387").Value = '$/unidad'
$ws.Range("O387").Value = 'Provincia de Quillota'
$ws.Range("P387").Value = 585
$ws.Range("Q387").Value = 1
$ws.Range("R387").Value = 'Hortaliza'
